$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Advance the weekly timesheet dates by one day each (new week: 0201-0207)
$ws.Range("B5").Value = 44228
$ws.Range("C5").Value = 44229
$ws.Range("D5").Value = 44230
$ws.Range("E5").Value = 44231
$ws.Range("F5").Value = 44232
$ws.Range("G5").Value = 44233
$ws.Range("H5").Value = 44234

# Move the active selection from E9 to B6
$ws.Range("B6").Select() | Out-Null
